$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'" + '63.138.84'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'" + '  +0.24%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'" + '2.548.70'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'" + '  +4.94%  '
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'" + '  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'" + '568.46'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'" + '  +0.70%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'" + '147.89'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'" + '  +4.95%  '
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'" + '  +0.04%  '
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'" + '  -1.22%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'" + '2.548.50'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'" + '  +4.98%  '
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'" + '  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'" + '  -1.60%  '
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'" + '  +0.36%  '
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'" + '  +0.84%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'" + '27.49'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'" + '  +4.98%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'" + '3.007.40'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'" + '  +5.04%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'" + '63.111.27'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'" + '  +0.35%  '
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'" + '  -0.45%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'" + '2.548.79'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'" + '  +4.82%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'" + '11.51'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'" + '  +2.61%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'" + '335.79'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'" + '  -1.18%  '
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'" + '  +1.26%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'" + '6.76'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'" + '  -0.37%  '
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'" + '  +0.10%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'" + '65.23'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'" + '  +0.13%  '
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'" + '  -2.49%  '
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'" + '  +4.80%  '
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'" + '  -0.20%  '
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'" + '8.40'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'" + '  +3.60%  '
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'" + '  +10.35%  '
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'" + '7.18'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'" + '  +8.66%  '
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "'" + '  +2.48%  '
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'" + '  +2.14%  '
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'" + '177.54'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'" + '  +1.80%  '
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "'" + '  +7.63%  '
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'" + '413.32'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'" + '  +11.86%  '
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'" + '  +0.29%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'" + '18.85'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'" + '  +1.16%  '
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'" + '4.39'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'" + '  -1.59%  '
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'" + '  -0.01%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'" + '1.75'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'" + '  +3.80%  '
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'" + '  +0.10%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'" + '39.37'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'" + '  -1.46%  '
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'" + '152.27'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'" + '  +2.31%  '
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'" + '  +1.68%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'" + '20.68'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'" + '  +0.55%  '
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'" + '0.609'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'" + '  +2.90%  '
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'" + '0.0966'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'" + '  +1.02%  '
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'" + '  +0.40%  '
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'" + '0.0235'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'" + '  +5.04%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'" + '18.33'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'" + '  +3.18%  '
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'" + '  +3.18%  '
$c.Style = "Normal"

